# Adds a new "RU" (Russia) region block at the top of the location list and
# a new "RER" (Europe) block at the bottom, by shifting every existing
# region block's location down one slot and re-appending a fresh copy of
# the original "RER" block's 24 rows at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate the still-untouched final block (rows 253:276, the "RER"
#    activity) down to rows 278:301 *before* any of the location values are
#    rewritten, so the appended copy keeps the original RER content/format.
$src = $ws.Range("A253:J276")
$dst = $ws.Range("A278:J301")
$src.Copy($dst)

# 2) Shift every block's "location" cells down one slot in the region
#    sequence (JP -> RAF -> RME -> RAS -> CA -> RLA -> GLO -> US -> IN ->
#    CN -> RER), and introduce the new first region, "RU", in the slot
#    that JP used to occupy.

# Block 1 (was JP) -> RU
$ws.Range("B8").Value = "RU"
$ws.Range("C16").Value = "RU"
$ws.Range("C20").Value = "RU"
$ws.Range("C22").Value = "RU"

# Block 2 (was RAF) -> JP
$ws.Range("B33").Value = "JP"
$ws.Range("C41").Value = "JP"
$ws.Range("C45").Value = "CN"
$ws.Range("C46").Value = "RoW"
$ws.Range("A47").Value = "market for electricity, medium voltage"
$ws.Range("C47").Value = "JP"

# Block 3 (was RME) -> RAF
$ws.Range("B58").Value = "RAF"
$ws.Range("C66").Value = "RAF"
$ws.Range("C70").Value = "RoW"
$ws.Range("C71").Value = "ZA"
$ws.Range("C72").Value = "RAF"

# Block 4 (was RAS) -> RME
$ws.Range("B83").Value = "RME"
$ws.Range("C91").Value = "RME"
$ws.Range("C95").Value = "RU"
$ws.Range("C96").Value = "RoW"
$ws.Range("C97").Value = "RME"

# Block 5 (was CA) -> RAS
$ws.Range("B108").Value = "RAS"
$ws.Range("C116").Value = "RAS"
$ws.Range("C120").Value = "CN"
$ws.Range("C121").Value = "CN"
$ws.Range("C122").Value = "RAS"

# Block 6 (was RLA) -> CA
$ws.Range("B133").Value = "CA"
$ws.Range("C141").Value = "CA"
$ws.Range("C145").Value = "RNA"
$ws.Range("C146").Value = "US"
$ws.Range("C147").Value = "CA"

# Block 7 (was GLO) -> RLA
$ws.Range("B158").Value = "RLA"
$ws.Range("C166").Value = "RLA"
$ws.Range("C170").Value = "RLA"
$ws.Range("C172").Value = "RLA"

# Block 8 (was US) -> GLO
$ws.Range("B183").Value = "GLO"
$ws.Range("C191").Value = "GLO"
$ws.Range("C195").Value = "RoW"
$ws.Range("C196").Value = "RoW"
$ws.Range("C197").Value = "GLO"

# Block 9 (was IN) -> US
$ws.Range("B208").Value = "US"
$ws.Range("C216").Value = "US"
$ws.Range("C220").Value = "RNA"
$ws.Range("C221").Value = "US"
$ws.Range("C222").Value = "US"

# Block 10 (was CN) -> IN
$ws.Range("B233").Value = "IN"
$ws.Range("C241").Value = "IN"
$ws.Range("C245").Value = "IN"
$ws.Range("C246").Value = "IN"
$ws.Range("C247").Value = "IN"

# Block 11 (was RER) -> CN
$ws.Range("B258").Value = "CN"
$ws.Range("C266").Value = "CN"
$ws.Range("C267").Value = "RoW"

# C270 previously carried a stray explicit font style (s="3"); the CN
# version of this row drops that formatting back to the sheet default.
$c270 = $ws.Range("C270")
$c270.Style = "Normal"
$c270.Value = "CN"

$ws.Range("C271").Value = "CN"
$ws.Range("C272").Value = "CN"
$ws.Range("C273").Value = "RoW"
$ws.Range("C274").Value = "RoW"
$ws.Range("C275").Value = "RoW"
$ws.Range("C276").Value = "RoW"

# 3) Restore the selection/scroll state to match the edited workbook: the
#    previous save had scrolled to the bottom (topLeftCell A258, selection
#    A276); the new save is scrolled back to the top with C23 selected.
$ws.Range("C23").Select()
